$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.473899126052856
$ws.Range("B1").Value = 1.502867698669434
$ws.Range("C1").Value = 8.067152976989746
$ws.Range("D1").Value = 2.100384712219238
$ws.Range("E1").Value = 1.082147002220154
